# "update template to calculate fsos for all regions"
#
# The three Facings-SOS KPI rows (rows 2-4 of the "Functional KPIs" sheet)
# currently carry a store attribute filter of address_city = Tokyo
# (columns O/P, "store_attr_1_name" / "store_attr_1_value"). Removing that
# filter makes those KPIs calculate for every region instead of being
# scoped to Tokyo only.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Functional KPIs")

# Row 2 - CCJP_FSOS_MANUF_BY_ALL_MANUF_IN_SCENE_TYPE
$ws.Range("O2").ClearContents()
$ws.Range("P2").ClearContents()

# Row 3 - CCJP_FSOS_MANUF_CAT_BY_ALL_MANUF_CAT_IN_SCENE_TYPE
$ws.Range("O3").ClearContents()
$ws.Range("P3").ClearContents()

# Row 4 - Facings SOS
$ws.Range("O4").ClearContents()
$ws.Range("P4").ClearContents()

# Leave the selection where the author left it after making the edit.
$ws.Range("S16").Select()
